$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '24.167.11'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  -2.24%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = "'" + '1.645.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -2.20%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('D4').Value = "'" + '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.17%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = "'" + '308.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -1.73%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = "'" + '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  -0.18%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = "'" + '0.3913'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.55%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = "'" + '0.3861'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -2.65%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = "'" + '1.002'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.14%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').Value = "'" + '1.373'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -3.61%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = "'" + '49.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -4.39%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = "'" + '0.08600'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -1.08%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = "'" + '23.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -6.59%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = "'" + '7.092'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  -3.44%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = "'" + '0.00001289'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -3.13%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').Value = "'" + '7.501'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -3.91%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').Value = "'" + '1.647.86'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  -4.07%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = "'" + '94.75'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +0.98%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('D19').Value = "'" + '0.06909'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -2.73%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('D20').Value = "'" + '20.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +0.28%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('D21').Value = "'" + '6.911'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -3.16%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('D22').Value = "'" + '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -0.26%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('E23').Value = "'" + '  -3.24%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('D24').Value = "'" + '24.172.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  -2.23%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('D25').Value = "'" + '2.425'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +2.94%  '
$ws.Range('E25').Style = 'Normal'

# Row 26
$ws.Range('D26').Value = "'" + '2.852'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +2.37%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('D27').Value = "'" + '22.33'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -5.28%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('D28').Value = "'" + '158.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -2.12%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('B29').Value = 'HuobiToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D29').Value = "'" + '5.386'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -7.26%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'" + '140.42'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -7.31%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'" + '8.193'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +4.17%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('D32').Value = "'" + '2.467'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +4.00%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('D33').Value = "'" + '1.831.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +3.23%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('D34').Value = "'" + '0.08170'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -3.40%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('D35').Value = "'" + '6.851'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -1.44%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('D36').Value = "'" + '0.02912'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -5.73%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = "'" + '0.9632'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -4.95%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = "'" + '0.2692'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  -3.91%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('D39').Value = "'" + '0.09196'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -3.56%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('D40').Value = "'" + '10.32'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -2.34%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('D41').Value = "'" + '1.458'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -1.61%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('D42').Value = "'" + '0.7557'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  -5.21%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('D43').Value = "'" + '13.04'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -4.85%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('D44').Value = "'" + '16.15'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -2.85%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('D45').Value = "'" + '0.6920'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -3.51%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('D46').Value = "'" + '2.465'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  -4.65%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('D47').Value = "'" + '4.096'
$ws.Range('D47').Style = 'Normal'

# Row 48
$ws.Range('D48').Value = "'" + '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  -0.15%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('D49').Value = "'" + '0.08369'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -3.58%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('D50').Value = "'" + '1.267'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  -4.93%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').Value = "'" + '133.68'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -3.61%  '
$ws.Range('E51').Style = 'Normal'
